# CIERRE 21 JUL 23
# Advance the payroll workbook from "SEMANA 28" (10-16 Jul 2023) to
# "SEMANA 29" (17-23 Jul 2023): update the week label, zero out last
# week's now-settled K4 figure, and record the new EXTRAS amount.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week banner (B9) drives H9 / B28 / H28 / B46 / H46 through =B9-style
# formulas, so updating this single cell ripples everywhere it is echoed.
$ws.Range("B9").Value = "SEMANA  29        DEL    17     Al   23  DE   JULIO    2023"

# Bonus/extra column for the first employee block resets to 0 for the
# new week (K7 = SUM(K4:K6) recalculates automatically).
$ws.Range("K4").Value = 0

# EXTRAS for the second employee block (E26 = SUM(E23:E25) recalculates
# automatically to 4200).
$ws.Range("E25").Value = 1400

# Move the on-screen selection to where the editor left off.
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("E26").Select() | Out-Null
